# MasterGallerySubmission.xlsx - add a "not_ok" sample row to the
# "Relative Samples" sheet and move the active selection there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relative Samples")

# New row of sample data (row 6), mirroring the existing rows 2-3 layout
# but with an extra 4th column and a text marker in column B.
$ws.Range("A6").Value = 2000
$ws.Range("B6").Value = "not_ok"
$ws.Range("C6").Value = 1000
$ws.Range("D6").Value = 2000

# Move the selection/active cell to D11, as recorded in the saved view state.
$ws.Activate() | Out-Null
$ws.Range("D11").Select() | Out-Null
